$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MD row (row 2) counters
$ws.Range("B2").Value = 9578

# C2 retains the same numeric value (2287) but needs its stored
# representation normalized (no more explicit numeric type / decimal).
# Force a change then set it back so the cell is rewritten.
$ws.Range("C2").Value = 0
$ws.Range("C2").Value = 2287

# Move the active selection to D11
$ws.Range("D11").Select()
